$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 8831.355113164813
$ws.Range("D5").Value = 8831.355113164813
$ws.Range("D9").Value = 11540.04331251616
$ws.Range("D10").Value = 11540.04331251616
$ws.Range("D14").Value = 11386.24488683539
$ws.Range("D15").Value = 11386.24488683539
